$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values that would otherwise be auto-parsed as numbers by Excel
# (single-dot decimals) are entered with a leading apostrophe so they stay
# text, matching the sheet's existing inlineStr/text storage for that column.
$ws.Range("D2").Value = "66.540.36"
$ws.Range("E2").Value = "  -4.50%  "
$ws.Range("D3").Value = "3.308.09"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'571.22"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").Value = "'181.49"
$ws.Range("E6").Value = "  -6.28%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -4.89%  "
$ws.Range("D12").Value = "3.885.10"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'27.10"
$ws.Range("E14").Value = "  -4.76%  "
$ws.Range("D15").Value = "66.604.37"
$ws.Range("E15").Value = "  -4.39%  "
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "3.279.92"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "'430.66"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'9.04"
$ws.Range("E27").Value = "  -5.70%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D30").Value = "'22.74"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").Value = "'6.77"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'159.96"
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "'27.12"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "2.813.23"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'0.788"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").Value = "'6.19"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "'0.0674"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'40.13"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "'24.35"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("E46").Value = "  -6.98%  "
$ws.Range("D47").Value = "'320.22"
$ws.Range("E47").Value = "  -7.35%  "
$ws.Range("E48").Value = "  -4.37%  "
$ws.Range("D49").Value = "'0.983"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").Value = "'6.17"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -1.30%  "
